$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New checks to append (Jenkins API / Pipeline rows).
# The underlying shared-string table must gain the four new unique strings
# in this exact order: DB_PORT..., IMAGESECRET_NAME..., The HELIX...UPDATE,
# HELIX...UPGRADE. That order does not match the row order they are used
# in, so the text is written into the worksheet's cells in the order that
# reproduces that shared-string sequence, while the row/column positions
# are set to their final home.

$ws.Cells.Item(138, 2).Value = "DB_PORT not blank"
$ws.Cells.Item(141, 2).Value = "IMAGESECRET_NAME not blank"
$ws.Cells.Item(140, 2).Value = "The HELIX_FULL_STACK_UPGRADE option sbould not be selected when the DEPLOYMENT_MODE is UPDATE"
$ws.Cells.Item(139, 2).Value = "HELIX_FULL_STACK_UPGRADE is not selected but it is required when the DEPLOYMENT_MODE is UPGRADE"

$newChecks = @(
    @{ Row = 138; Num = 240 },
    @{ Row = 139; Num = 241 },
    @{ Row = 140; Num = 242 },
    @{ Row = 141; Num = 243 }
)

foreach ($item in $newChecks) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Num
    $ws.Cells.Item($r, 3).Value = "Jenkins API"
    $ws.Cells.Item($r, 4).Value = "Pipeline"
}

$ws.Range("D141").Select() | Out-Null
